$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 2177.1304
$ws.Range("I28").Value = 1239.6471
$ws.Range("J28").Value = 4833.3335
$ws.Range("K28").Value = 1239.6471
$ws.Range("L28").Value = 4833.3335
$ws.Range("M28").Value = -754.6470999999999
$ws.Range("N28").Value = -5803.3335

# Row 74
$ws.Range("H74").Value = 4500
$ws.Range("I74").Value = 4000
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 4000
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -3064
$ws.Range("N74").Value = -6872

# Row 77
$ws.Range("H77").Value = 4500
$ws.Range("I77").Value = 4000
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 20000
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -15320
$ws.Range("N77").Value = -34360

# Row 132
$ws.Range("H132").Value = 32681732
$ws.Range("I132").Value = 4831995.5
$ws.Range("J132").Value = 90913000
$ws.Range("K132").Value = 14495986.5
$ws.Range("L132").Value = 272739000
$ws.Range("M132").Value = -14493456.5
$ws.Range("N132").Value = -272744060

# Row 135
$ws.Range("H135").Value = 454.30768
$ws.Range("I135").Value = 454.30768
$ws.Range("K135").Value = 4088.76912
$ws.Range("M135").Value = -1553.76912

# Row 137
$ws.Range("H137").Value = 568250.8
$ws.Range("I137").Value = 1454.2632
$ws.Range("J137").Value = 927221.9399999999
$ws.Range("K137").Value = 4362.7896
$ws.Range("L137").Value = 2781665.82
$ws.Range("M137").Value = -1812.7896
$ws.Range("N137").Value = -2786765.82

# Row 138
$ws.Range("H138").Value = 4003861
$ws.Range("I138").Value = 2296.0667
$ws.Range("J138").Value = 5718817.5
$ws.Range("K138").Value = 6888.2001
$ws.Range("L138").Value = 17156452.5
$ws.Range("M138").Value = -1748.2001
$ws.Range("N138").Value = -17166732.5

# Row 141
$ws.Range("H141").Value = 4746
$ws.Range("I141").Value = 2942.1428
$ws.Range("K141").Value = 8826.428400000001
$ws.Range("M141").Value = -3646.428400000001


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3487.2678
$ws.Range("I32").Value = 2909.3655
$ws.Range("J32").Value = 11000
$ws.Range("K32").Value = 2909.3655
$ws.Range("L32").Value = 11000
$ws.Range("M32").Value = -2622.3655
$ws.Range("N32").Value = -11574

# Row 41
$ws.Range("H41").Value = 8528
$ws.Range("I41").Value = 2056
$ws.Range("J41").Value = 15000
$ws.Range("K41").Value = 2056
$ws.Range("L41").Value = 15000
$ws.Range("M41").Value = -1642
$ws.Range("N41").Value = -15828

# Row 61
$ws.Range("H61").Value = 3077.077
$ws.Range("I61").Value = 1141
$ws.Range("J61").Value = 4287.125
$ws.Range("K61").Value = 1141
$ws.Range("L61").Value = 4287.125
$ws.Range("M61").Value = -929
$ws.Range("N61").Value = -4711.125

# Row 63
$ws.Range("H63").Value = 3237.1875
$ws.Range("I63").Value = 2179.5
$ws.Range("K63").Value = 2179.5
$ws.Range("M63").Value = -1493.5

# Row 66
$ws.Range("H66").Value = 3237.1875
$ws.Range("I66").Value = 2179.5
$ws.Range("K66").Value = 10897.5
$ws.Range("M66").Value = -7465.5

# Row 74
$ws.Range("H74").Value = 84851.836
$ws.Range("I74").Value = 251106
$ws.Range("J74").Value = 1724.75
$ws.Range("K74").Value = 251106
$ws.Range("L74").Value = 1724.75
$ws.Range("M74").Value = -250232
$ws.Range("N74").Value = -3472.75

# Row 77
$ws.Range("H77").Value = 84851.836
$ws.Range("I77").Value = 251106
$ws.Range("J77").Value = 1724.75
$ws.Range("K77").Value = 1255530
$ws.Range("L77").Value = 8623.75
$ws.Range("M77").Value = -1251162
$ws.Range("N77").Value = -17359.75

# Row 108
$ws.Range("H108").Value = 422896.8
$ws.Range("J108").Value = 422896.8
$ws.Range("L108").Value = 422896.8
$ws.Range("N108").Value = -430576.8

# Row 132
$ws.Range("H132").Value = 2458
$ws.Range("I132").Value = 1898.1578
$ws.Range("J132").Value = 3977.5715
$ws.Range("K132").Value = 5694.4734
$ws.Range("L132").Value = 11932.7145
$ws.Range("M132").Value = -3164.4734
$ws.Range("N132").Value = -16992.7145

# Row 136
$ws.Range("H136").Value = 3077.077
$ws.Range("I136").Value = 1141
$ws.Range("J136").Value = 4287.125
$ws.Range("K136").Value = 3423
$ws.Range("L136").Value = 12861.375
$ws.Range("M136").Value = -873
$ws.Range("N136").Value = -17961.375


$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1909.762
$ws.Range("I86").Value = 1712.0588
$ws.Range("K86").Value = 1712.0588
$ws.Range("M86").Value = -589.0588

# Row 89
$ws.Range("H89").Value = 1909.762
$ws.Range("I89").Value = 1712.0588
$ws.Range("K89").Value = 8560.294
$ws.Range("M89").Value = -2944.294

# Row 134
$ws.Range("H134").Value = 4178.9575
$ws.Range("I134").Value = 3982.1282
$ws.Range("J134").Value = 5138.5
$ws.Range("K134").Value = 11946.3846
$ws.Range("L134").Value = 15415.5
$ws.Range("M134").Value = -9411.384600000001
$ws.Range("N134").Value = -20485.5


$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 475.75
$ws.Range("I22").Value = 360.63635
$ws.Range("J22").Value = 729
$ws.Range("K22").Value = 360.63635
$ws.Range("L22").Value = 729
$ws.Range("M22").Value = -10.63634999999999
$ws.Range("N22").Value = -1429

# Row 44
$ws.Range("H44").Value = 12383.2
$ws.Range("I44").Value = 11489
$ws.Range("J44").Value = 12979.333
$ws.Range("K44").Value = 11489
$ws.Range("L44").Value = 12979.333
$ws.Range("M44").Value = -11047
$ws.Range("N44").Value = -13863.333

# Row 55
$ws.Range("H55").Value = 50000
$ws.Range("J55").Value = 50000
$ws.Range("L55").Value = 50000
$ws.Range("N55").Value = -50630

# Row 134
$ws.Range("H134").Value = 33335330
$ws.Range("I134").Value = 4002011.2
$ws.Range("K134").Value = 12006033.6
$ws.Range("M134").Value = -12003498.6


$ws = $wb.Worksheets.Item("CUL")
# Row 20
$ws.Range("H20").Value = 1798.625
$ws.Range("I20").Value = 598.1667
$ws.Range("K20").Value = 1794.5001
$ws.Range("M20").Value = -1567.5001

# Row 33
$ws.Range("H33").Value = 78.92308
$ws.Range("I33").Value = 27.714285
$ws.Range("J33").Value = 138.66667
$ws.Range("K33").Value = 166.28571
$ws.Range("L33").Value = 832.0000200000001
$ws.Range("M33").Value = 116.71429
$ws.Range("N33").Value = -1398.00002

# Row 97
$ws.Range("H97").Value = 713.931
$ws.Range("I97").Value = 457.14285
$ws.Range("J97").Value = 795.63635
$ws.Range("K97").Value = 1371.42855
$ws.Range("L97").Value = 2386.90905
$ws.Range("M97").Value = -875.4285500000001
$ws.Range("N97").Value = -3378.90905

# Row 100
$ws.Range("H100").Value = 3979.7334
$ws.Range("J100").Value = 3979.7334
$ws.Range("L100").Value = 11939.2002
$ws.Range("N100").Value = -13561.2002

# Row 104
$ws.Range("H104").Value = 2500
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 2500
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 7500
$ws.Range("M104").ClearContents()
$ws.Range("N104").Value = -12742

# Row 106
$ws.Range("H106").Value = 4085.1853
$ws.Range("J106").Value = 4085.1853
$ws.Range("L106").Value = 12255.5559
$ws.Range("N106").Value = -14147.5559

# Row 113
$ws.Range("H113").Value = 562.1111
$ws.Range("I113").Value = 560.86664
$ws.Range("J113").Value = 563.6667
$ws.Range("K113").Value = 1682.59992
$ws.Range("L113").Value = 1691.0001
$ws.Range("M113").Value = 487.4000800000001
$ws.Range("N113").Value = -6031.0001

# Row 115
$ws.Range("H115").Value = 3105.6
$ws.Range("I115").Value = 1028
$ws.Range("J115").Value = 3625
$ws.Range("K115").Value = 3084
$ws.Range("L115").Value = 10875
$ws.Range("M115").Value = -1909
$ws.Range("N115").Value = -13225

# Row 119
$ws.Range("H119").Value = 6776.3335
$ws.Range("I119").Value = 6776.3335
$ws.Range("K119").Value = 20329.0005
$ws.Range("M119").Value = -15491.0005


$ws = $wb.Worksheets.Item("LTW")
# Row 26
$ws.Range("H26").Value = 11340
$ws.Range("I26").Value = 10000
$ws.Range("K26").Value = 10000
$ws.Range("M26").Value = -9705

# Row 132
$ws.Range("H132").Value = 3227.1455
$ws.Range("I132").Value = 2892.5806
$ws.Range("J132").Value = 3659.2917
$ws.Range("K132").Value = 8677.7418
$ws.Range("L132").Value = 10977.8751
$ws.Range("M132").Value = -6147.7418
$ws.Range("N132").Value = -16037.8751

# Row 136
$ws.Range("H136").Value = 1595.8334
$ws.Range("I136").Value = 1584.8529
$ws.Range("J136").Value = 1614.5
$ws.Range("K136").Value = 4754.5587
$ws.Range("L136").Value = 4843.5
$ws.Range("M136").Value = -2204.5587
$ws.Range("N136").Value = -9943.5


$ws = $wb.Worksheets.Item("WVR")
# Row 10
$ws.Range("H10").Value = 13000
$ws.Range("I10").Value = 1000
$ws.Range("J10").Value = 19000
$ws.Range("K10").Value = 1000
$ws.Range("L10").Value = 19000
$ws.Range("M10").Value = -831
$ws.Range("N10").Value = -19338

# Row 11
$ws.Range("H11").Value = 51000
$ws.Range("I11").Value = 64000
$ws.Range("J11").Value = 25000
$ws.Range("K11").Value = 64000
$ws.Range("L11").Value = 25000
$ws.Range("M11").Value = -63858
$ws.Range("N11").Value = -25284

# Row 13
$ws.Range("H13").Value = 1305
$ws.Range("I13").Value = 1305
$ws.Range("K13").Value = 1305
$ws.Range("M13").Value = -1165

# Row 100
$ws.Range("H100").Value = 940.26666
$ws.Range("I100").Value = 924.4
$ws.Range("J100").Value = 972
$ws.Range("K100").Value = 1848.8
$ws.Range("L100").Value = 1944
$ws.Range("M100").Value = -1307.8
$ws.Range("N100").Value = -3026

# Row 107
$ws.Range("H107").Value = 23973
$ws.Range("I107").Value = 38178.8
$ws.Range("J107").Value = 296.66666
$ws.Range("K107").Value = 114536.4
$ws.Range("L107").Value = 889.9999799999999
$ws.Range("M107").Value = -112616.4
$ws.Range("N107").Value = -4729.99998

# Row 136
$ws.Range("H136").Value = 3514.5518
$ws.Range("I136").Value = 3883.7646
$ws.Range("J136").Value = 2991.5
$ws.Range("K136").Value = 11651.2938
$ws.Range("L136").Value = 8974.5
$ws.Range("M136").Value = -9101.293799999999
$ws.Range("N136").Value = -14074.5

